$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I7").Value = "test_none_to_true_correctly_categorises_nutrients()"
$ws.Range("I8").Value = "test_false_to_true_correctly_categorises_nutrients()"
$ws.Range("I9").Value = "test_true_to_false_with_single_nutrient_correctly_categorises_opposing_implication()"

[void]$ws.Range("I9").Select()
